$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 155, shifting existing rows 155-183 down to 156-184
$ws.Rows.Item(155).Insert()

# Populate the newly inserted row 155 with the new record.
# Columns A,B,C,E,F,G,I,R carry the same constant values used throughout this block.
$ws.Range("A155").Value = 11
$ws.Range("B155").Value = "Vega Monumental Concepción"
$ws.Range("C155").Value = "Bíobío"
$ws.Range("D155").Value = 45015
$ws.Range("D155").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E155").Value = 8
$ws.Range("F155").Value = 100112021
$ws.Range("G155").Value = "Ají"
$ws.Range("H155").Value = "Americana (o)"
$ws.Range("I155").Value = "Primera"
$ws.Range("J155").Value = 50
$ws.Range("K155").Value = 6500
$ws.Range("L155").Value = 7000
$ws.Range("M155").Value = 6700
$ws.Range("N155").Value = "$/saco 25 kilos"
$ws.Range("O155").Value = "Región Metropolitana"
$ws.Range("P155").Value = 268
$ws.Range("Q155").Value = 25
$ws.Range("R155").Value = "Hortaliza"
